$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2, shifting the existing rows (and their
# formatting) down by one.
$ws.Rows(2).Insert()

# The new row 2 holds the latest date; prices are unchanged from the
# previous day, so just duplicate them.
$cellA = $ws.Range("A2")

# Force text so the ISO-formatted date string isn't auto-converted into a
# date serial number, then drop the temporary number format so the cell's
# style stays the default (matching the rest of the column).
$cellA.NumberFormat = "@"
$cellA.Value = "2026-01-26"
$cellA.ClearFormats()

$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610
